$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G and H columns (and I2) with new computed return values
$ws.Cells.Item(2, 7).Value = 0.1025027763869769
$ws.Cells.Item(2, 8).Value = 22.21955496760163
$ws.Cells.Item(2, 9).Value = -51.5775013732102
$ws.Cells.Item(3, 7).Value = 0.1133155111762229
$ws.Cells.Item(3, 8).Value = -2.612857540317413
$ws.Cells.Item(4, 7).Value = -0.06715218242078227
$ws.Cells.Item(4, 8).Value = -335.3385338325562
$ws.Cells.Item(5, 7).Value = -0.06749119842334067
$ws.Cells.Item(5, 8).Value = 6.068718184744195
$ws.Cells.Item(6, 7).Value = 0.01804675787153686
$ws.Cells.Item(6, 8).Value = -49.22937060532719
$ws.Cells.Item(7, 7).Value = -0.006239031931259804
$ws.Cells.Item(7, 8).Value = -130.7415703674701
$ws.Cells.Item(8, 7).Value = -0.1447150755211115
$ws.Cells.Item(8, 8).Value = -2.563595106666171
$ws.Cells.Item(9, 7).Value = -0.1498241635018284
$ws.Cells.Item(9, 8).Value = -9.450972201822946
$ws.Cells.Item(10, 7).Value = -0.07296410981027383
$ws.Cells.Item(10, 8).Value = 32.45760167340618
$ws.Cells.Item(11, 7).Value = -0.1125935480041143
$ws.Cells.Item(11, 8).Value = -69.20963334993955
$ws.Cells.Item(12, 7).Value = -0.3441437782566001
$ws.Cells.Item(12, 8).Value = 17.11299959828504
$ws.Cells.Item(13, 7).Value = -0.4213463864745443
$ws.Cells.Item(13, 8).Value = 6.137901341738974
$ws.Cells.Item(14, 7).Value = -0.06646568423815459
$ws.Cells.Item(14, 8).Value = -31.0037891688375
$ws.Cells.Item(15, 7).Value = 0.05104570234941084
$ws.Cells.Item(15, 8).Value = 161.7303362349937
$ws.Cells.Item(16, 7).Value = 0.09843664764422394
$ws.Cells.Item(16, 8).Value = -32.58811320440682
$ws.Cells.Item(17, 7).Value = 0.1613556916008601
$ws.Cells.Item(17, 8).Value = 31.58843348462021
$ws.Cells.Item(18, 7).Value = 0.1381801339839694
$ws.Cells.Item(18, 8).Value = -0.2000955523147199
$ws.Cells.Item(19, 7).Value = 0.1280312699007649
$ws.Cells.Item(19, 8).Value = 34.19519778926346
$ws.Cells.Item(20, 7).Value = 0.01988877917143014
$ws.Cells.Item(20, 8).Value = -22.41266008438647
$ws.Cells.Item(21, 7).Value = 0.03348184890156319
$ws.Cells.Item(21, 8).Value = -55.45678452163874
$ws.Cells.Item(24, 7).Value = 0.1035953615945692
$ws.Cells.Item(24, 8).Value = 3.14512592821406
$ws.Cells.Item(25, 7).Value = 0.1533948793181145
$ws.Cells.Item(25, 8).Value = 1.219655287329703
$ws.Cells.Item(26, 7).Value = 0.07817252341080842
$ws.Cells.Item(26, 8).Value = -1.2034016429184
$ws.Cells.Item(27, 7).Value = 0.0772204763637018
$ws.Cells.Item(27, 8).Value = -22.66554895517652
$ws.Cells.Item(28, 7).Value = -0.2317431113365898
$ws.Cells.Item(28, 8).Value = -8.738737809338845
$ws.Cells.Item(29, 7).Value = -0.1754012841147151
$ws.Cells.Item(29, 8).Value = 14.56620521939482
$ws.Cells.Item(30, 7).Value = 0.05034708199228426
$ws.Cells.Item(30, 8).Value = 14.08386128683146
$ws.Cells.Item(31, 7).Value = 0.02985552001031537
$ws.Cells.Item(31, 8).Value = 13.37634601386778
$ws.Cells.Item(32, 7).Value = 0.09921783541826959
$ws.Cells.Item(32, 8).Value = 4.483260156933307
$ws.Cells.Item(33, 7).Value = 0.1162509996657865
$ws.Cells.Item(33, 8).Value = 11.81755982226482
$ws.Cells.Item(34, 7).Value = 0.02134662738435218
$ws.Cells.Item(34, 8).Value = -54.0217771210587
$ws.Cells.Item(35, 7).Value = 0.02457999633784065
$ws.Cells.Item(35, 8).Value = 224.4229786340618
$ws.Cells.Item(36, 7).Value = 0.05281910749423492
$ws.Cells.Item(36, 8).Value = -8.516553779363056
$ws.Cells.Item(37, 7).Value = 0.06939669136618949
$ws.Cells.Item(37, 8).Value = -1.32165329741461
$ws.Cells.Item(38, 7).Value = 0.01717896915461914
$ws.Cells.Item(38, 8).Value = -67.20720372732801
$ws.Cells.Item(39, 7).Value = 0.03137824639659204
$ws.Cells.Item(39, 8).Value = 51.31286147379132
$ws.Cells.Item(40, 7).Value = 0.006677234343228602
$ws.Cells.Item(40, 8).Value = 178.6390599925196
$ws.Cells.Item(41, 7).Value = 0.03427881320986012
$ws.Cells.Item(41, 8).Value = -3.046287047291473
$ws.Cells.Item(42, 7).Value = 0.1315523975460724
$ws.Cells.Item(42, 8).Value = -1.596073060346122
$ws.Cells.Item(43, 7).Value = 0.1527579874918305
$ws.Cells.Item(43, 8).Value = 2.538425323795694
$ws.Cells.Item(44, 7).Value = -0.008772972536095085
$ws.Cells.Item(44, 8).Value = -3.072218060791526
$ws.Cells.Item(45, 7).Value = -0.0007497534914795126
$ws.Cells.Item(45, 8).Value = 93.17156005654302
$ws.Cells.Item(46, 7).Value = -0.006803958036114187
$ws.Cells.Item(46, 8).Value = -106.615948417758
$ws.Cells.Item(47, 7).Value = -0.007415655403798909
$ws.Cells.Item(47, 8).Value = 20.07683878775171
$ws.Cells.Item(48, 7).Value = 0.06015623641561887
$ws.Cells.Item(48, 8).Value = 19.66124757345275
$ws.Cells.Item(49, 7).Value = 0.06206414010345319
$ws.Cells.Item(49, 8).Value = -6.055637073058213
$ws.Cells.Item(50, 7).Value = 0.1482925717581902
$ws.Cells.Item(50, 8).Value = -8.03999445745206
$ws.Cells.Item(51, 7).Value = 0.1548138082436349
$ws.Cells.Item(51, 8).Value = -9.52643256814895
$ws.Cells.Item(52, 7).Value = -0.1715577974739761
$ws.Cells.Item(52, 8).Value = -6.935564662968637
$ws.Cells.Item(53, 7).Value = -0.1294462584653656
$ws.Cells.Item(53, 8).Value = -2.690707671391977
$ws.Cells.Item(54, 7).Value = 0.106692625326148
$ws.Cells.Item(54, 8).Value = 13.83898704179775
$ws.Cells.Item(55, 7).Value = 0.1112856976080556
$ws.Cells.Item(55, 8).Value = -1.587404654150349
$ws.Cells.Item(56, 7).Value = -0.02256687437287529
$ws.Cells.Item(56, 8).Value = -209.1025977103067
$ws.Cells.Item(57, 7).Value = -0.02637834686491473
$ws.Cells.Item(57, 8).Value = -15.35354683308103
$ws.Cells.Item(58, 7).Value = 0.04962168846540198
$ws.Cells.Item(58, 8).Value = -11.99951654555907
$ws.Cells.Item(59, 7).Value = 0.07047688978439125
$ws.Cells.Item(59, 8).Value = -1.870004970945611
$ws.Cells.Item(60, 7).Value = 0.07526740431833416
$ws.Cells.Item(60, 8).Value = 7.564887523670885
$ws.Cells.Item(61, 7).Value = 0.07002927229736566
$ws.Cells.Item(61, 8).Value = 47.34695712348514
$ws.Cells.Item(62, 7).Value = 0.05672765506653829
$ws.Cells.Item(62, 8).Value = -22.2486288081259
$ws.Cells.Item(63, 7).Value = 0.06020291762802309
$ws.Cells.Item(63, 8).Value = -7.937465579875576
$ws.Cells.Item(64, 7).Value = -0.0299662225131561
$ws.Cells.Item(64, 8).Value = 27.64492268412403
$ws.Cells.Item(65, 7).Value = 0.009694747136775605
$ws.Cells.Item(65, 8).Value = 119.6553495078447
$ws.Cells.Item(66, 7).Value = 0.04419971712923735
$ws.Cells.Item(66, 8).Value = 133.4412936442983
$ws.Cells.Item(67, 7).Value = -0.00574132886836934
$ws.Cells.Item(67, 8).Value = -121.9554643577433
$ws.Cells.Item(68, 7).Value = -0.02588825344596269
$ws.Cells.Item(68, 8).Value = -4642.170381750907
$ws.Cells.Item(69, 7).Value = 0.01391485300303564
$ws.Cells.Item(69, 8).Value = 207.7642688223956
$ws.Cells.Item(70, 7).Value = -0.04338650056215061
$ws.Cells.Item(70, 8).Value = -58.06759287256104
$ws.Cells.Item(71, 7).Value = -0.05618025430817563
$ws.Cells.Item(71, 8).Value = -1.975430709621432
$ws.Cells.Item(72, 7).Value = -0.1467221920145771
$ws.Cells.Item(72, 8).Value = 1.086704716849592
$ws.Cells.Item(73, 7).Value = -0.1594948385759187
$ws.Cells.Item(73, 8).Value = -10.156101407872
$ws.Cells.Item(74, 7).Value = 0.1272297770201104
$ws.Cells.Item(74, 8).Value = 0.9350742132193048
$ws.Cells.Item(75, 7).Value = 0.1508622974937688
$ws.Cells.Item(75, 8).Value = 11.6076723898914
$ws.Cells.Item(76, 7).Value = -0.04146851843065622
$ws.Cells.Item(76, 8).Value = -20.41519428984922
$ws.Cells.Item(77, 7).Value = -0.06310320470778977
$ws.Cells.Item(77, 8).Value = -36.6254443820475
$ws.Cells.Item(78, 7).Value = 0.08610082721727157
$ws.Cells.Item(78, 8).Value = -6.586318379301447
$ws.Cells.Item(79, 7).Value = 0.08707362625885948
$ws.Cells.Item(79, 8).Value = -9.772404505491501
$ws.Cells.Item(80, 7).Value = -0.1719484673256274
$ws.Cells.Item(80, 8).Value = -5.871215510629816
$ws.Cells.Item(81, 7).Value = -0.1955303726566369
$ws.Cells.Item(81, 8).Value = 9.664233948859858
$ws.Cells.Item(82, 7).Value = 0.162597476642189
$ws.Cells.Item(82, 8).Value = 17.19051329725603
$ws.Cells.Item(83, 7).Value = 0.1908089206885961
$ws.Cells.Item(83, 8).Value = 15.90899656458023
$ws.Cells.Item(84, 7).Value = 0.05743076112559856
$ws.Cells.Item(84, 8).Value = 310.3654735876539
$ws.Cells.Item(85, 7).Value = 0.04943497850950154
$ws.Cells.Item(85, 8).Value = 118.3602516549894
